$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Delete the "Still need to add to add outdoor location example..." paragraph
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Still need to add to add outdoor location example*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2. Remove the _GoBack bookmark from the "WALKING CALCULATION." paragraph
#    (it gets re-added at the very end of the document later on)
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ---------------------------------------------------------------------
# 3. Replace the last 3 trailing empty paragraphs with new content.
#    (the first 2 of the 5 trailing empty paragraphs are left untouched)
# ---------------------------------------------------------------------
$total = $d.Paragraphs.Count
$delStart = $d.Paragraphs.Item($total - 2).Range.Start
$delEnd = $d.Paragraphs.Item($total).Range.End
$d.Range($delStart, $delEnd).Delete()

# ---------------------------------------------------------------------
# 4. Append the new simple (single or multi run) paragraphs
# ---------------------------------------------------------------------

# "Still TODO"
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Still TODO"

# "Future applications- reference Colin, speak to Derek and possibly recruitment."
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Future applications- reference Colin, speak to Derek and possibly recruitment."

# "Future Work"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Future Work"

# "Route Modifiers" + "- Look over"  (2 runs)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Route Modifiers"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "- Look over"
$mergePos = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.End - 1
$d.Range($mergePos, $mergePos + 1).Delete()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

# "Location Services"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Location Services"

# "Disjoint Routes" + "-Flesh out" (2 runs)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Disjoint Routes"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "-Flesh out"
$mergePos = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.End - 1
$d.Range($mergePos, $mergePos + 1).Delete()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

# "Reaching more platforms" + " " + "Flesh out" (3 runs)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Reaching more platforms"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = " "
$mergePos = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.End - 1
$d.Range($mergePos, $mergePos + 1).Delete()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Flesh out"
$mergePos = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.End - 1
$d.Range($mergePos, $mergePos + 1).Delete()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

# "Release and Updates" + " " + "Flesh out" (3 runs)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Release and Updates"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = " "
$mergePos = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.End - 1
$d.Range($mergePos, $mergePos + 1).Delete()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Flesh out"
$mergePos = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.End - 1
$d.Range($mergePos, $mergePos + 1).Delete()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

# "Performance Analysis & Future Proofing- Add results and tweak"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Performance Analysis & Future Proofing- Add results and tweak"

# "Testing Strategies " + en-dash + " Add tests and write up"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Testing Strategies " + [char]0x2013 + " Add tests and write up"

# "Make it super clear throughout this actually works lots of screen shots of it working and looking pretty if possible."
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Make it super clear throughout this actually works lots of screen shots of it working and looking pretty if possible."

# "Abstract"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Abstract"

# "Pimp out introduction speaking about what happens in each chapter and why you<rsquo>d be interested"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "Pimp out introduction speaking about what happens in each chapter and why you" + [char]0x2019 + "d be interested"

Write-Host "Checkpoint2 paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Host "$i`: [$($pp.Range.Text)]"
}
